# New weekly price observation added for "Femacal de La Calera - Poroto verde".
# A new row is inserted at row 363 (pushing the existing rows 363-412 down to
# 364-413, growing the used range from A1:R412 to A1:R413), and the new row
# is populated with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 363, shifting rows 363:412 down to 364:413.
$ws.Rows.Item(363).Insert()

# Populate the newly inserted row 363 with the new observation.
$ws.Range("A363").Value = 3
$ws.Range("B363").Value = "Femacal de La Calera"
$ws.Range("C363").Value = "Coquimbo"
$ws.Range("D363").Value = 44776
$ws.Range("E363").Value = 5
$ws.Range("F363").Value = 100112031
$ws.Range("G363").Value = "Poroto verde"
$ws.Range("H363").Value = "Magnum"
$ws.Range("I363").Value = "Primera"
$ws.Range("J363").Value = 83
$ws.Range("K363").Value = 32000
$ws.Range("L363").Value = 33000
$ws.Range("M363").Value = 32542
$ws.Range("N363").Value = "`$/malla 25 kilos"
$ws.Range("O363").Value = "Región de Arica y Parinacota"
$ws.Range("P363").Value = 1302
$ws.Range("Q363").Value = 25
$ws.Range("R363").Value = "Hortaliza"
